$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# survey sheet: add a new "over 2 years" branch (if/else/end if) driven
# by the ADA custom-date field, with Norwegian + Portuguese labels.
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Row 4 gains a styled (shaded) blank cell in column A, matching the
# existing "begin screen" / "end screen" marker rows above it.
$survey.Range("A2").Copy() | Out-Null
$survey.Range("A4").PasteSpecial(-4122) | Out-Null
$survey.Range("B4").Value = "end screen"

# Row 5: new "begin screen" marker row.
$survey.Range("A2").Copy() | Out-Null
$survey.Range("A5").PasteSpecial(-4122) | Out-Null
$survey.Range("B5").Value = "begin screen"

# Row 6: clause "if" with the condition testing the ADA date.
$survey.Range("B6").Value = "if"
$survey.Range("C6").Value = "data('ADA') && data('ADA').length>0 && moment().diff(moment(data('ADA'), '\D:DD,\M:MM,\Y:YYYY'),'years')>=2"

# Row 7: note prompt shown when the branch condition is true.
$survey.Range("D7").Value = "note"
$survey.Range("F7").Value = "lblKnown"
$survey.Range("G7").Value = "OVER 2 ÅR"
$survey.Range("H7").Value = "MAS DE 2 ANOS"

# Row 8: clause "else".
$survey.Range("B8").Value = "else"

# Row 9: note prompt shown in the else branch.
$survey.Range("D9").Value = "note"
$survey.Range("F9").Value = "lblUnknown"
$survey.Range("G9").Value = "IKKE OVER 2 ÅR"
$survey.Range("H9").Value = "NAO MAS DE 2 ANOS"

# Row 10: clause "end if".
$survey.Range("B10").Value = "end if"

# Row 11: closing "end screen" marker row.
$survey.Range("B11").Value = "end screen"

# ---------------------------------------------------------------------
# model sheet: register the new "TT" session variable (string type).
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$model.Range("A3").Value = "TT"
$model.Range("B3").Value = "string"
$model.Range("C3").Value = $false

# ---------------------------------------------------------------------
# Leave the survey sheet as the active tab/selection, as in the edit.
# ---------------------------------------------------------------------
$model.Range("D3").Select() | Out-Null
$survey.Activate()
$survey.Range("C7").Select() | Out-Null
